$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2999.5935
$ws.Range("I32").Value = 2989.9155
$ws.Range("J32").Value = 3033.95
$ws.Range("K32").Value = 2989.9155
$ws.Range("L32").Value = 3033.95
$ws.Range("M32").Value = -2702.9155
$ws.Range("N32").Value = -3607.95

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 49000
$ws.Range("J88").Value = 49000
$ws.Range("L88").Value = 49000
$ws.Range("N88").Value = -49812

$ws.Range("H91").Value = 49000
$ws.Range("J91").Value = 49000
$ws.Range("L91").Value = 49000
$ws.Range("N91").Value = -51808

$ws.Range("H105").Value = 1825.4429
$ws.Range("I105").Value = 1827.9701
$ws.Range("K105").Value = 1827.9701
$ws.Range("M105").Value = -80.9701

$ws.Range("H107").Value = 1076.375
$ws.Range("I107").Value = 796.2778
$ws.Range("J107").Value = 1916.6666
$ws.Range("K107").Value = 796.2778
$ws.Range("L107").Value = 1916.6666
$ws.Range("M107").Value = 1123.7222
$ws.Range("N107").Value = -5756.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 680.8461
$ws.Range("I22").Value = 421.1579
$ws.Range("J22").Value = 1385.7142
$ws.Range("K22").Value = 421.1579
$ws.Range("L22").Value = 1385.7142
$ws.Range("M22").Value = -71.15789999999998
$ws.Range("N22").Value = -2085.7142

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5682404.5
$ws.Range("I113").Value = 626.1667
$ws.Range("J113").Value = 12500539
$ws.Range("K113").Value = 1878.5001
$ws.Range("L113").Value = 37501617
$ws.Range("M113").Value = 291.4999
$ws.Range("N113").Value = -37505957

$ws.Range("H132").Value = 1971.1818
$ws.Range("I132").Value = 845.3077
$ws.Range("J132").Value = 2703
$ws.Range("K132").Value = 7607.7693
$ws.Range("L132").Value = 24327
$ws.Range("M132").Value = -5077.7693
$ws.Range("N132").Value = -29387

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1716.9166
$ws.Range("I61").Value = 1720.8
$ws.Range("J61").Value = 1714.1428
$ws.Range("K61").Value = 1720.8
$ws.Range("L61").Value = 1714.1428
$ws.Range("M61").Value = -1518.8
$ws.Range("N61").Value = -2118.1428

$ws.Range("H68").Value = 1182.5574
$ws.Range("I68").Value = 973.2593000000001
$ws.Range("J68").Value = 2797.1428
$ws.Range("K68").Value = 973.2593000000001
$ws.Range("L68").Value = 2797.1428
$ws.Range("M68").Value = -224.2593000000001
$ws.Range("N68").Value = -4295.1428

$ws.Range("H71").Value = 1182.5574
$ws.Range("I71").Value = 973.2593000000001
$ws.Range("J71").Value = 2797.1428
$ws.Range("K71").Value = 4866.2965
$ws.Range("L71").Value = 13985.714
$ws.Range("M71").Value = -1122.2965
$ws.Range("N71").Value = -21473.714

$ws.Range("H92").Value = 29000
$ws.Range("J92").Value = 29000
$ws.Range("L92").Value = 29000
$ws.Range("N92").Value = -33992

$ws.Range("H93").Value = 1966
$ws.Range("I93").Value = 1564.2667
$ws.Range("J93").Value = 2970.3333
$ws.Range("K93").Value = 1564.2667
$ws.Range("L93").Value = 2970.3333
$ws.Range("M93").Value = -316.2666999999999
$ws.Range("N93").Value = -5466.3333

$ws.Range("H94").Value = 29307.273
$ws.Range("J94").Value = 29307.273
$ws.Range("L94").Value = 29307.273
$ws.Range("N94").Value = -30659.273

$ws.Range("H95").Value = 33872
$ws.Range("J95").Value = 33872
$ws.Range("L95").Value = 33872
$ws.Range("N95").Value = -39364

$ws.Range("H96").Value = 38766.668
$ws.Range("J96").Value = 38766.668
$ws.Range("L96").Value = 38766.668
$ws.Range("N96").Value = -44258.668

$ws.Range("H97").Value = 34475
$ws.Range("J97").Value = 34475
$ws.Range("L97").Value = 34475
$ws.Range("N97").Value = -36457

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H99").Value = 22419.666
$ws.Range("I99").Value = 16129.5
$ws.Range("K99").Value = 16129.5
$ws.Range("M99").Value = -13134.5

$ws.Range("H100").Value = 2450.3
$ws.Range("I100").Value = 2099.8333
$ws.Range("K100").Value = 2099.8333
$ws.Range("M100").Value = -1558.8333

$ws.Range("H101").Value = 35362
$ws.Range("J101").Value = 35362
$ws.Range("L101").Value = 35362
$ws.Range("N101").Value = -41852

$ws.Range("H104").Value = 23925.715
$ws.Range("J104").Value = 23925.715
$ws.Range("L104").Value = 23925.715
$ws.Range("N104").Value = -30913.715

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H106").Value = 29956.5
$ws.Range("J106").Value = 29956.5
$ws.Range("L106").Value = 29956.5
$ws.Range("N106").Value = -32480.5

$ws.Range("H108").Value = 49950
$ws.Range("J108").Value = 49950
$ws.Range("L108").Value = 49950
$ws.Range("N108").Value = -57630

$ws.Range("H110").Value = 31571.285
$ws.Range("J110").Value = 31571.285
$ws.Range("L110").Value = 31571.285
$ws.Range("N110").Value = -39751.285

$ws.Range("H111").Value = 31304.334
$ws.Range("J111").Value = 31304.334
$ws.Range("L111").Value = 31304.334
$ws.Range("N111").Value = -39484.334

$ws.Range("H113").Value = 1716.9166
$ws.Range("I113").Value = 1720.8
$ws.Range("J113").Value = 1714.1428
$ws.Range("K113").Value = 1720.8
$ws.Range("L113").Value = 1714.1428
$ws.Range("M113").Value = 449.2
$ws.Range("N113").Value = -6054.1428

$ws.Range("H114").Value = 39715
$ws.Range("J114").Value = 39715
$ws.Range("L114").Value = 39715
$ws.Range("N114").Value = -48393

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 39900
$ws.Range("J80").Value = 39900
$ws.Range("L80").Value = 39900
$ws.Range("N80").Value = -41896

$ws.Range("H83").Value = 39900
$ws.Range("J83").Value = 39900
$ws.Range("L83").Value = 119700
$ws.Range("N83").Value = -129684

$ws.Range("H118").Value = 29318.572
$ws.Range("J118").Value = 29318.572
$ws.Range("L118").Value = 29318.572
$ws.Range("N118").Value = -29318.572
